$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -----------------------------------------------------------------
# 1. Update existing qqq_return / qld_return values (columns B and C)
# -----------------------------------------------------------------
$ws.Range("B2").Value = 0.172725464013983
$ws.Range("C2").Value = 0.32380718115198

$ws.Range("B3").Value = 0.312276012629918
$ws.Range("C3").Value = 0.593909113633052

$ws.Range("B4").Value = 0.51889597120594
$ws.Range("C4").Value = 1.06809994620302

$ws.Range("B5").Value = 1.03148385529191
$ws.Range("C5").Value = 2.29321491486981

$ws.Range("B6").Value = 2.69324169105462
$ws.Range("C6").Value = 4.62879679790476

$ws.Range("B7").Value = 4.8228797429133
$ws.Range("C7").Value = 9.33449855580181

$ws.Range("B8").Value = 7.79636739815915
$ws.Range("C8").Value = 12.2610603510839

$ws.Range("B9").Value = 16.2575243960344
$ws.Range("C9").Value = 40.5537579139876

$ws.Range("B10").Value = 39.1502108778063
$ws.Range("C10").Value = 179.093558135726

# -----------------------------------------------------------------
# 2. Add the new tqqq_return column (D), reusing the existing header /
#    data cell formatting from column C so no new style entries are
#    introduced.
# -----------------------------------------------------------------
$ws.Range("C1").Copy()
$ws.Range("D1").PasteSpecial(-4122)
$ws.Range("D1").Value = "tqqq_return"

$ws.Range("C2:C10").Copy()
$ws.Range("D2:D10").PasteSpecial(-4122)

$ws.Range("D2").Value = 0.438057065078531
$ws.Range("D3").Value = 0.797561756343892
$ws.Range("D4").Value = 1.51658129584604
$ws.Range("D5").Value = 3.28179609241294
$ws.Range("D6").Value = 4.83020106621927
$ws.Range("D7").Value = 2.91128496484833
$ws.Range("D8").Value = 2.19263625022717
$ws.Range("D9").Value = 12.5325406536171
$ws.Range("D10").Value = 84.8793314679889

$wb.Application.CutCopyMode = $false

# -----------------------------------------------------------------
# 3. Percent formatting now shows two decimal places (0.00%) across
#    the header and data cells of B:D.
# -----------------------------------------------------------------
$ws.Range("B1:D1").NumberFormat = "0.00%"
$ws.Range("B2:D10").NumberFormat = "0.00%"

# -----------------------------------------------------------------
# 4. Column widths were tightened slightly (target character widths:
#    A=13.6363636363636, B:D=12.8181818181818 - the engine quantizes
#    stored widths to 1/6 character increments, so these ColumnWidth
#    inputs are the closest achievable approximations).
# -----------------------------------------------------------------
$ws.Range("A1").ColumnWidth = 12.833333333333332
$ws.Range("B1").ColumnWidth = 12
$ws.Range("C1").ColumnWidth = 12
$ws.Range("D1").ColumnWidth = 12

# -----------------------------------------------------------------
# 5. The built-in "Hyperlink" cell style was renamed to "Link".
# -----------------------------------------------------------------
$hyperlinkStyle = $wb.Styles.Item("Hyperlink")
$hyperlinkStyle.Name = "Link"

# -----------------------------------------------------------------
# 6. Selection moved.
# -----------------------------------------------------------------
$ws.Range("J5").Select()
